$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the text of the "PROPOSTA DO MVP" and "FUNCIONALIDADES" blocks.
#    PROPOSTA DO MVP loses its trailing blank lines.
#    FUNCIONALIDADES loses the "Upload de imagens e videos" line.
#    Both keep their header bold, the remaining text regular (matches the
#    formatting already used by the other blocks on the sheet).
#    (PROPOSTA DO MVP is rewritten before FUNCIONALIDADES so the new shared
#    strings end up appended in that same relative order.)
# ---------------------------------------------------------------------------

$propostaHeader = "PROPOSTA DO MVP"
$propostaBody = "Rede social voltada para contatos entre influenciadores digitais e empresas que desejam acompanhar e/ou contratar os serviços de influenciadores digitais."
$propostaText = $propostaHeader + [char]10 + $propostaBody

$cellProposta = $ws.Range("B2")
$cellProposta.Value = $propostaText
$headerLen2 = $propostaHeader.Length
$totalLen2 = $propostaText.Length
$cellProposta.Characters(1, $headerLen2).Font.Bold = $true
$bodyRun2 = $cellProposta.Characters($headerLen2 + 1, $totalLen2 - $headerLen2)
$bodyRun2.Font.Name = "Calibri"
$bodyRun2.Font.Size = 11
$bodyRun2.Font.Bold = $false

$funcionalidadesHeader = "FUNCIONALIDADES"
$funcionalidadesBody = "Cadastro de influenciadores digitais" + [char]10 + "Cadastro de empresas" + [char]10 + "Envio de notificações"
$funcionalidadesText = $funcionalidadesHeader + [char]10 + $funcionalidadesBody

$cellFuncionalidades = $ws.Range("B4")
$cellFuncionalidades.Value = $funcionalidadesText
$headerLen = $funcionalidadesHeader.Length
$totalLen = $funcionalidadesText.Length
$cellFuncionalidades.Characters(1, $headerLen).Font.Bold = $true
$bodyRun = $cellFuncionalidades.Characters($headerLen + 1, $totalLen - $headerLen)
$bodyRun.Font.Name = "Calibri"
$bodyRun.Font.Size = 11
$bodyRun.Font.Bold = $false

# ---------------------------------------------------------------------------
# 2. Vertically center every cell of the canvas table (title, headers and
#    blank filler cells alike now share a consistent vertical="center").
# ---------------------------------------------------------------------------

$ws.Range("A1:C7").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3. Resize columns / rows.
# ---------------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 28.666666666666668
$ws.Columns.Item(2).ColumnWidth = 60.166666666666664
$ws.Columns.Item(3).ColumnWidth = 28.666666666666668

$ws.Rows.Item(1).RowHeight = 36
$ws.Rows.Item(2).RowHeight = 80.4
$ws.Rows.Item(4).RowHeight = 50.4
$ws.Rows.Item(7).RowHeight = 95.4

# ---------------------------------------------------------------------------
# 4. Update the selection shown when the sheet is reopened.
# ---------------------------------------------------------------------------

$ws.Range("B4:B5").Select()
